# Slide 6 ("STGraph - Operations"), Content Placeholder 2: update wording per the
# commit diff (valid-path formula gains the "e" subscript, the ">GraphNode-TSNode"
# bullet becomes an intro line, and the following bullets get trailing punctuation).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$tr.Paragraphs(2, 1).Runs(1).Text = "Path(ni, …, nk) è valido ⇔ ⋂j=i..k-1 Ie(nj,nj+1) ≠ ∅"
$tr.Paragraphs(4, 1).Runs(1).Text = "Traversing a virtual edge:"
$tr.Paragraphs(5, 1).Runs(1).Text = "Entails a query to AsterixDB ;"
$tr.Paragraphs(6, 1).Runs(1).Text = "Filter pushdown ;"
$tr.Paragraphs(7, 1).Runs(1).Text = "No support for cross time-series operations ."
